$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 (SMH): fill in purchase details ---
# Copy number formats from an existing filled row so the date/percent
# formatting matches the rest of the table.
$ws.Range("C2").Copy()
$ws.Range("C4").PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Range("E4").PasteSpecial(-4122)
$ws.Range("F2").Copy()
$ws.Range("F4").PasteSpecial(-4122)

$ws.Range("B4").Value = 255.6
$ws.Range("C4").Value = 44384
$ws.Range("D4").Value = 2
$ws.Range("E4").Value = 0.2
$ws.Range("F4").Value = 0.05

# --- Row 6 (ACES): clear out the purchase details ---
$ws.Range("B6:C6").ClearContents()
$ws.Range("E6:F6").ClearContents()
$ws.Range("D6").Value = 0

# --- Row 7 (ARKW): update share count ---
$ws.Range("D7").Value = 2

# --- Row 10 (GXTG): fill in purchase details ---
$ws.Range("C2").Copy()
$ws.Range("C10").PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Range("E10").PasteSpecial(-4122)
$ws.Range("F2").Copy()
$ws.Range("F10").PasteSpecial(-4122)

$ws.Range("B10").Value = 51.25
$ws.Range("C10").Value = 44384
$ws.Range("D10").Value = 6
$ws.Range("E10").Value = 0.2
$ws.Range("F10").Value = 0.05

# --- Row 11 (FANG): clear out the purchase details ---
$ws.Range("B11:C11").ClearContents()
$ws.Range("E11:F11").ClearContents()
$ws.Range("D11").Value = 0

# --- Row 12 (LIT): update share count ---
$ws.Range("D12").Value = 4

# --- Update active selection ---
$ws.Range("D2").Select()
